# Generate Report for Handoff
# Updates the localization-status report: status flips from the
# "handed back" state to "ready for handoff", and the handoff/generation
# timestamps are refreshed to the new run's datetime. Also tightens the
# "Status" column widths that used to be sized for the long status text.

$wb = $excel.ActiveWorkbook

# ---- Overview sheet --------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsOverview.Range("G2").Value = "2016-08-24 11:02:16"
$wsOverview.Columns.Item(5).ColumnWidth = 17
$wsOverview.Columns.Item(6).ColumnWidth = 17

# ---- zh-cn sheet -------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("H2").Value = "2016-08-24 11:02:03"
$wsZhCn.Columns.Item(3).ColumnWidth = 17

# ---- de-de sheet -------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("H2").Value = "2016-08-24 11:02:16"
$wsDeDe.Columns.Item(3).ColumnWidth = 17
